# Add a new forecast-origin date column (2020-05-07, column AB) and a new
# forecast-target date row (2020-05-21, row 40) to both the "cases" and
# "deaths" sheets, and populate the diagonal of new forecast / observed
# values that came in with that day's update.

$wb = $excel.ActiveWorkbook

function Set-ForecastSheet {
    param($ws, $Observed26, $AAUpdates, $ABValues, $AB40)

    # New column header: AB1 = the new forecast-origin date string.
    # Force text formatting first so Excel doesn't auto-convert the
    # ISO-looking string into a date serial number (matches the existing
    # header cells C1:AA1, which are stored as plain text).
    $ws.Range("AB1").NumberFormat = "@"
    $ws.Range("AB1").Value = "2020-05-07"
    $ws.Range("AB1").Style = "Normal"

    # Newly-observed value for the 2020-05-07 row (row 26), column B
    # ("Observed").
    $ws.Range("B26").Value = $Observed26

    # Existing AA column (forecast made on 2020-05-06) got revised values
    # for every still-open target row.
    foreach ($row in $AAUpdates.Keys) {
        $ws.Range("AA$row").Value = $AAUpdates[$row]
    }

    # New AB column (forecast made on 2020-05-07) values for every
    # still-open target row (27 through the new row 40).
    foreach ($row in $ABValues.Keys) {
        $ws.Range("AB$row").Value = $ABValues[$row]
    }

    # New target-date row: 2020-05-21 (row 40). Column A is the row label
    # (again forced to text so it is not reinterpreted as a date), and the
    # only populated data cell is the brand-new AB column.
    $ws.Range("A40").NumberFormat = "@"
    $ws.Range("A40").Value = "2020-05-21"
    $ws.Range("A40").Style = "Normal"
    $ws.Range("AB40").Value = $AB40
}

# ---- cases sheet ----
$wsCases = $wb.Worksheets.Item("cases")

$casesAA = @{
    26 = 14285
    27 = 15200
    28 = 16177
    29 = 17097
    30 = 18062
    31 = 19122
    32 = 20076
    33 = 21024
    34 = 21628
    35 = 22508
    36 = 23046
    37 = 23663
    38 = 24188
    39 = 24839
}
$casesAB = @{
    27 = 15161
    28 = 15857
    29 = 17021
    30 = 17826
    31 = 18625
    32 = 19422
    33 = 20175
    34 = 20656
    35 = 21376
    36 = 21914
    37 = 22352
    38 = 22855
    39 = 23328
}

Set-ForecastSheet $wsCases 14156 $casesAA $casesAB 23990

# ---- deaths sheet ----
$wsDeaths = $wb.Worksheets.Item("deaths")

$deathsAA = @{}
$deathsAB = @{
    27 = 1465
    28 = 1555
    29 = 1664
    30 = 1753
    31 = 1859
    32 = 1962
    33 = 2065
    34 = 2136
    35 = 2221
    36 = 2281
    37 = 2343
    38 = 2409
    39 = 2484
}

Set-ForecastSheet $wsDeaths 1394 $deathsAA $deathsAB 2568
